$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Computer Vision")

# --- Row 30: new paper entry (Byte Latent Transformer) ---

# A30: title, rich text with bold prefix "ByteLatentTransformer" then regular
# ": PatchesScaleBetter ThanTokens"
$titleCell = $ws.Range("A30")
$titleCell.Value2 = "ByteLatentTransformer: PatchesScaleBetter ThanTokens"
$titleBold = $titleCell.Characters(1, 21)
$titleBold.Font.Bold = $true
$titleRest = $titleCell.Characters(22, 31)
$titleRest.Font.Bold = $false

# B30: author
$ws.Range("B30").Value2 = "Artidoro Pagnoni et al."

# C30: year
$ws.Range("C30").Value2 = 2024

# D30: conference / venue
$ws.Range("D30").Value2 = "Preprint"

# E30: paper link (text + hyperlink)
$linkCell = $ws.Range("E30")
$linkCell.Value2 = "https://ai.meta.com/research/publications/byte-latent-transformer-patches-scale-better-than-tokens/"
$ws.Hyperlinks.Add($linkCell, "https://ai.meta.com/research/publications/byte-latent-transformer-patches-scale-better-than-tokens/") | Out-Null

# F30: what?
$ws.Range("F30").Value2 = "This paper introduces byte latent transformer (BLT) that improves the scaling trends of LLMs. Operating in the byte space is costly due to long sequence lengths. Previous methods have explored more efficient attention mechanisms and remove attention all together. BLT is dynamic and learnable method for grouping bytes  into patchs (a patch-based approach). One unique difference between this method and other patch methods is that BLT has no fixed size or vocabulary for patches. "

# G30: contributions?
$ws.Range("G30").Value2 = "(1) BLT, byte latent LLM architecture that dynamically allocates compute to improve FLOP (floating point operations per seconds) (2) Unlocking a new dimension of scaling LLMs . (3) Authors demonstrate imporved robustness of BLT models to input noise and awareness of sub-word aspects"

# H30: category
$ws.Range("H30").Value2 = "LLMs/ Architecture"

# Row height: the new row wraps to a tall multi-line row like its neighbours
$ws.Rows.Item(30).RowHeight = 187

# --- View state tweaks matching the author's final scroll/zoom position ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 38
$ws.Range("H31").Select()
